$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a D-column (Price) cell as text so Excel does not auto-convert
# the punctuation-as-thousands-separator strings (e.g. "67.700.28") or
# trailing-zero / scientific-notation-prone values into numbers.
function Set-PriceText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-PriceText "D2" "67.700.28"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.787.13"
$ws.Range("E3").Value = "  +1.52%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-PriceText "D5" "595.09"
$ws.Range("E5").Value = "  +0.36%  "

# Row 6 - Solana
Set-PriceText "D6" "166.14"
$ws.Range("E6").Value = "  -0.60%  "

# Row 7 - LidoStakedEther
Set-PriceText "D7" "3.786.29"
$ws.Range("E7").Value = "  +1.43%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.14%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.16%  "

# Row 11 - Toncoin
Set-PriceText "D11" "6.36"
$ws.Range("E11").Value = "  -1.57%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.17%  "

# Row 13 - ShibaInu
Set-PriceText "D13" "0.0000256"
$ws.Range("E13").Value = "  -1.16%  "

# Row 14 - Avalanche
Set-PriceText "D14" "36.19"
$ws.Range("E14").Value = "  -0.08%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText "D15" "4.419.39"
$ws.Range("E15").Value = "  +1.46%  "

# Row 16 - WrappedEther
Set-PriceText "D16" "3.784.26"
$ws.Range("E16").Value = "  +1.50%  "

# Row 17 - now Chainlink (was WrappedBTC)
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-PriceText "D17" "18.43"
$ws.Range("E17").Value = "  +3.02%  "

# Row 18 - now WrappedBTC (was Chainlink)
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-PriceText "D18" "67.673.51"
$ws.Range("E18").Value = "  -0.73%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.42%  "

# Row 20 - Polkadot
Set-PriceText "D20" "6.97"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21 - Uniswap
Set-PriceText "D21" "10.19"
$ws.Range("E21").Value = "  -4.88%  "

# Row 22 - BitcoinCash
Set-PriceText "D22" "458.05"
$ws.Range("E22").Value = "  -1.90%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.696"
$ws.Range("E23").Value = "  -0.56%  "

# Row 24 - PEPE
Set-PriceText "D24" "0.0000152"
$ws.Range("E24").Value = "  +5.34%  "

# Row 25 - Litecoin
Set-PriceText "D25" "83.65"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26 - InternetComputer(DFINITY)
Set-PriceText "D26" "11.87"
$ws.Range("E26").Value = "  -1.66%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  -2.73%  "

# Row 28 - RenderToken
Set-PriceText "D28" "10.07"
$ws.Range("E28").Value = "  -0.34%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.13%  "

# Row 30 - PancakeSwap
Set-PriceText "D30" "2.78"
$ws.Range("E30").Value = "  +0.14%  "

# Row 31 - NEARProtocol
Set-PriceText "D31" "7.30"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32 - EthereumClassic
Set-PriceText "D32" "29.90"
$ws.Range("E32").Value = "  +0.15%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  +0.46%  "

# Row 34 - Aptos
Set-PriceText "D34" "9.17"
$ws.Range("E34").Value = "  -0.95%  "

# Row 35 - Binance-PegBSC-USD
Set-PriceText "D35" "0.999"
$ws.Range("E35").Value = "  -0.32%  "

# Row 36 - RenzoRestakedETH
Set-PriceText "D36" "3.739.27"
$ws.Range("E36").Value = "  +1.49%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -1.12%  "

# Row 38 - dogwifhat
Set-PriceText "D38" "3.32"
$ws.Range("E38").Value = "  -1.98%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +0.71%  "

# Row 40 - Mantle
Set-PriceText "D40" "0.996"
$ws.Range("E40").Value = "  +0.07%  "

# Row 41 - Filecoin
Set-PriceText "D41" "5.75"
$ws.Range("E41").Value = "  -0.85%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.02%  "

# Row 44 - Arweave
Set-PriceText "D44" "44.53"
$ws.Range("E44").Value = "  +3.02%  "

# Row 45 - now TheGraph (was OKB)
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-PriceText "D45" "0.298"
$ws.Range("E45").Value = "  -2.45%  "

# Row 46 - now OKB (was TheGraph)
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-PriceText "D46" "47.15"
$ws.Range("E46").Value = "  +2.96%  "

# Row 47 - Cosmos
Set-PriceText "D47" "8.38"
$ws.Range("E47").Value = "  -2.48%  "

# Row 48 - Monero
Set-PriceText "D48" "147.67"
$ws.Range("E48").Value = "  +0.95%  "

# Row 49 - Bittensor
Set-PriceText "D49" "391.62"
$ws.Range("E49").Value = "  -0.44%  "

# Row 50 - Stacks
$ws.Range("E50").Value = "  -5.85%  "

# Row 51 - Maker
Set-PriceText "D51" "2.753.32"
$ws.Range("E51").Value = "  +2.67%  "
